$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 11690
$ws.Range("I21").Value = 7112.5
$ws.Range("J21").Value = 30000
$ws.Range("K21").Value = 7112.5
$ws.Range("L21").Value = 30000
$ws.Range("M21").Value = -6644.5
$ws.Range("N21").Value = -30936
$ws.Range("H23").Value = 11690
$ws.Range("I23").Value = 7112.5
$ws.Range("J23").Value = 30000
$ws.Range("K23").Value = 7112.5
$ws.Range("L23").Value = 30000
$ws.Range("M23").Value = -6878.5
$ws.Range("N23").Value = -30468
$ws.Range("H62").Value = 1697.5
$ws.Range("I62").Value = 1600
$ws.Range("J62").Value = 1990
$ws.Range("K62").Value = 1600
$ws.Range("L62").Value = 1990
$ws.Range("M62").Value = -976
$ws.Range("N62").Value = -3238
$ws.Range("H65").Value = 1697.5
$ws.Range("I65").Value = 1600
$ws.Range("J65").Value = 1990
$ws.Range("K65").Value = 8000
$ws.Range("L65").Value = 9950
$ws.Range("M65").Value = -4880
$ws.Range("N65").Value = -16190
$ws.Range("H76").Value = 3961.8572
$ws.Range("I76").Value = 3547
$ws.Range("K76").Value = 3547
$ws.Range("M76").Value = -3232
$ws.Range("H79").Value = 3961.8572
$ws.Range("I79").Value = 3547
$ws.Range("K79").Value = 3547
$ws.Range("M79").Value = -2455
$ws.Range("H93").Value = 44199.668
$ws.Range("J93").Value = 44199.668
$ws.Range("L93").Value = 44199.668
$ws.Range("N93").Value = -49191.668
$ws.Range("H94").Value = 3526
$ws.Range("I94").Value = 3526
$ws.Range("K94").Value = 3526
$ws.Range("M94").Value = -3075
$ws.Range("H103").Value = 589.6667
$ws.Range("I103").Value = 345
$ws.Range("J103").Value = 607.1429000000001
$ws.Range("K103").Value = 1035
$ws.Range("L103").Value = 1821.4287
$ws.Range("M103").Value = -449
$ws.Range("N103").Value = -2993.4287
$ws.Range("H108").Value = 94996.57000000001
$ws.Range("J108").Value = 94996.57000000001
$ws.Range("L108").Value = 94996.57000000001
$ws.Range("N108").Value = -102676.57
$ws.Range("H132").Value = 1326.5714
$ws.Range("I132").Value = 1255.4849
$ws.Range("K132").Value = 3766.4547
$ws.Range("M132").Value = -1236.4547
$ws.Range("H137").Value = 325213.47
$ws.Range("I137").Value = 1865.2084
$ws.Range("J137").Value = 694754.3
$ws.Range("K137").Value = 5595.6252
$ws.Range("L137").Value = 2084262.9
$ws.Range("M137").Value = -3045.6252
$ws.Range("N137").Value = -2089362.9

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 52724.9
$ws.Range("I61").Value = 2562.5625
$ws.Range("K61").Value = 2562.5625
$ws.Range("M61").Value = -2350.5625
$ws.Range("H63").Value = 5560346.5
$ws.Range("I63").Value = 8337018.5
$ws.Range("K63").Value = 8337018.5
$ws.Range("M63").Value = -8336332.5
$ws.Range("H66").Value = 5560346.5
$ws.Range("I66").Value = 8337018.5
$ws.Range("K66").Value = 41685092.5
$ws.Range("M66").Value = -41681660.5
$ws.Range("H74").Value = 3204.4866
$ws.Range("I74").Value = 1506.05
$ws.Range("K74").Value = 1506.05
$ws.Range("M74").Value = -632.05
$ws.Range("H77").Value = 3204.4866
$ws.Range("I77").Value = 1506.05
$ws.Range("K77").Value = 7530.25
$ws.Range("M77").Value = -3162.25
$ws.Range("H102").Value = 81007.86
$ws.Range("I102").Value = 93418.09
$ws.Range("J102").Value = 35503.668
$ws.Range("K102").Value = 93418.09
$ws.Range("L102").Value = 35503.668
$ws.Range("M102").Value = -91796.09
$ws.Range("N102").Value = -38747.668
$ws.Range("H105").Value = 112999.5
$ws.Range("J105").Value = 112999.5
$ws.Range("L105").Value = 112999.5
$ws.Range("N105").Value = -119987.5
$ws.Range("H132").Value = 2148.9167
$ws.Range("I132").Value = 1992
$ws.Range("K132").Value = 5976
$ws.Range("M132").Value = -3446
$ws.Range("H136").Value = 52724.9
$ws.Range("I136").Value = 2562.5625
$ws.Range("K136").Value = 7687.6875
$ws.Range("M136").Value = -5137.6875

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1584.7727
$ws.Range("I80").Value = 1549.5
$ws.Range("J80").Value = 1588.3
$ws.Range("K80").Value = 1549.5
$ws.Range("L80").Value = 1588.3
$ws.Range("M80").Value = -551.5
$ws.Range("N80").Value = -3584.3
$ws.Range("H83").Value = 1584.7727
$ws.Range("I83").Value = 1549.5
$ws.Range("J83").Value = 1588.3
$ws.Range("K83").Value = 7747.5
$ws.Range("L83").Value = 7941.5
$ws.Range("M83").Value = -2755.5
$ws.Range("N83").Value = -17925.5
$ws.Range("H94").Value = 1129.909
$ws.Range("I94").Value = 1183.1177
$ws.Range("J94").Value = 949
$ws.Range("K94").Value = 1183.1177
$ws.Range("L94").Value = 949
$ws.Range("M94").Value = -732.1177
$ws.Range("N94").Value = -1851
$ws.Range("H134").Value = 2995.4243
$ws.Range("I134").Value = 2220.9355
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 6662.806500000001
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = -4127.806500000001
$ws.Range("N134").Value = -50070

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4599.6665
$ws.Range("I58").Value = 4599.6665
$ws.Range("K58").Value = 4599.6665
$ws.Range("M58").Value = -4396.6665
$ws.Range("H94").Value = 1201.8
$ws.Range("I94").Value = 1249.5
$ws.Range("K94").Value = 1249.5
$ws.Range("M94").Value = -798.5
$ws.Range("H105").Value = 3154.7083
$ws.Range("I105").Value = 1541.9412
$ws.Range("K105").Value = 1541.9412
$ws.Range("M105").Value = 205.0588
$ws.Range("H132").Value = 918145.25
$ws.Range("I132").Value = 816434.6
$ws.Range("K132").Value = 2449303.8
$ws.Range("M132").Value = -2446773.8
$ws.Range("H136").Value = 4599.6665
$ws.Range("I136").Value = 4599.6665
$ws.Range("K136").Value = 13798.9995
$ws.Range("M136").Value = -11248.9995

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1457.375
$ws.Range("J68").Value = 1531.8
$ws.Range("L68").Value = 4595.4
$ws.Range("N68").Value = -6217.4
$ws.Range("H71").Value = 1457.375
$ws.Range("J71").Value = 1531.8
$ws.Range("L71").Value = 13786.2
$ws.Range("N71").Value = -21898.2

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 352000
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 352000
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H97").Value = 2095.0952
$ws.Range("I97").Value = 2337
$ws.Range("K97").Value = 2337
$ws.Range("M97").Value = -1841
$ws.Range("H132").Value = 4287.814
$ws.Range("I132").Value = 3414.9429
$ws.Range("K132").Value = 10244.8287
$ws.Range("M132").Value = -7714.8287

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 78055.7
$ws.Range("I22").Value = 3137.625
$ws.Range("K22").Value = 3137.625
$ws.Range("M22").Value = -2842.625
$ws.Range("H27").Value = 78055.7
$ws.Range("I27").Value = 3137.625
$ws.Range("K27").Value = 3137.625
$ws.Range("M27").Value = -3030.625
$ws.Range("H40").Value = 2782156.2
$ws.Range("I40").Value = 4455.154
$ws.Range("J40").Value = 7940744.5
$ws.Range("K40").Value = 4455.154
$ws.Range("L40").Value = 7940744.5
$ws.Range("M40").Value = -4319.154
$ws.Range("N40").Value = -7941016.5
$ws.Range("H68").Value = 3417.5
$ws.Range("J68").Value = 3301
$ws.Range("L68").Value = 3301
$ws.Range("N68").Value = -4799
$ws.Range("H71").Value = 3417.5
$ws.Range("J71").Value = 3301
$ws.Range("L71").Value = 16505
$ws.Range("N71").Value = -23993

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9599
$ws.Range("I62").Value = 9599
$ws.Range("K62").Value = 9599
$ws.Range("M62").Value = -8975
$ws.Range("H65").Value = 9599
$ws.Range("I65").Value = 9599
$ws.Range("K65").Value = 47995
$ws.Range("M65").Value = -44875
$ws.Range("H81").Value = 13984.143
$ws.Range("J81").Value = 10000
$ws.Range("L81").Value = 20000
$ws.Range("N81").Value = -22122
$ws.Range("H84").Value = 13984.143
$ws.Range("J84").Value = 10000
$ws.Range("L84").Value = 100000
$ws.Range("N84").Value = -110608
$ws.Range("H122").Value = 3704.0688
$ws.Range("I122").Value = 3799.5417
$ws.Range("J122").Value = 3245.8
$ws.Range("K122").Value = 11398.6251
$ws.Range("L122").Value = 9737.400000000001
$ws.Range("M122").Value = -8948.625100000001
$ws.Range("N122").Value = -14637.4
$ws.Range("H132").Value = 1567.7675
$ws.Range("I132").Value = 1152.0667
$ws.Range("J132").Value = 2527.077
$ws.Range("K132").Value = 3456.2001
$ws.Range("L132").Value = 7581.231000000001
$ws.Range("M132").Value = -926.2001
$ws.Range("N132").Value = -12641.231
